$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 39166.15
$ws.Range("I40").Value = 68599.92999999999
$ws.Range("J40").Value = 2373.9167
$ws.Range("K40").Value = 68599.92999999999
$ws.Range("L40").Value = 2373.9167
$ws.Range("M40").Value = -68424.92999999999
$ws.Range("N40").Value = -2723.9167
$ws.Range("H43").Value = 1956.2916
$ws.Range("I43").Value = 1847.3572
$ws.Range("J43").Value = 2108.8
$ws.Range("K43").Value = 1847.3572
$ws.Range("L43").Value = 2108.8
$ws.Range("M43").Value = -1778.3572
$ws.Range("N43").Value = -2246.8
$ws.Range("H52").Value = 203.6
$ws.Range("I52").Value = 203.6
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 610.8
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -450.8
$ws.Range("N52").ClearContents()
$ws.Range("H69").Value = 5904.5386
$ws.Range("I69").Value = 4199.75
$ws.Range("J69").Value = 6662.222
$ws.Range("K69").Value = 12599.25
$ws.Range("L69").Value = 19986.666
$ws.Range("M69").Value = -11725.25
$ws.Range("N69").Value = -21734.666
$ws.Range("H72").Value = 5904.5386
$ws.Range("I72").Value = 4199.75
$ws.Range("J72").Value = 6662.222
$ws.Range("K72").Value = 37797.75
$ws.Range("L72").Value = 59959.998
$ws.Range("M72").Value = -33429.75
$ws.Range("N72").Value = -68695.99799999999
$ws.Range("H132").Value = 4241632.5
$ws.Range("I132").Value = 4633820.5
$ws.Range("K132").Value = 13901461.5
$ws.Range("M132").Value = -13898931.5
$ws.Range("H138").Value = 5800.192
$ws.Range("I138").Value = 1738.0952
$ws.Range("J138").Value = 7440.654
$ws.Range("K138").Value = 5214.2856
$ws.Range("L138").Value = 22321.962
$ws.Range("M138").Value = -74.28560000000016
$ws.Range("N138").Value = -32601.962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 44002.715
$ws.Range("J23").Value = 29601.4
$ws.Range("L23").Value = 29601.4
$ws.Range("N23").Value = -30119.4
$ws.Range("H59").Value = 10750
$ws.Range("I59").Value = 5500
$ws.Range("J59").Value = 16000
$ws.Range("K59").Value = 5500
$ws.Range("L59").Value = 16000
$ws.Range("M59").Value = -4696
$ws.Range("N59").Value = -17608
$ws.Range("H74").Value = 3354.25
$ws.Range("I74").Value = 708.5
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 708.5
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = 165.5
$ws.Range("N74").Value = -7748
$ws.Range("H77").Value = 3354.25
$ws.Range("I77").Value = 708.5
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 3542.5
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = 825.5
$ws.Range("N77").Value = -38736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 113770.664
$ws.Range("J105").Value = 169283.67
$ws.Range("L105").Value = 169283.67
$ws.Range("N105").Value = -172777.67
$ws.Range("H106").Value = 16694.2
$ws.Range("J106").Value = 16694.2
$ws.Range("L106").Value = 16694.2
$ws.Range("N106").Value = -19218.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 17009.666
$ws.Range("I99").Value = 5865.4287
$ws.Range("K99").Value = 5865.4287
$ws.Range("M99").Value = -4367.4287
$ws.Range("H126").Value = 17009.666
$ws.Range("I126").Value = 5865.4287
$ws.Range("K126").Value = 17596.2861
$ws.Range("M126").Value = -15126.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 140
$ws.Range("I2").Value = 166.25
$ws.Range("J2").Value = 87.5
$ws.Range("K2").Value = 997.5
$ws.Range("L2").Value = 525
$ws.Range("M2").Value = -884.5
$ws.Range("N2").Value = -751
$ws.Range("H38").Value = 60.666668
$ws.Range("I38").Value = 15
$ws.Range("J38").Value = 106.333336
$ws.Range("K38").Value = 45
$ws.Range("L38").Value = 319.000008
$ws.Range("M38").Value = 302
$ws.Range("N38").Value = -1013.000008
$ws.Range("H58").Value = 760
$ws.Range("I58").Value = 576
$ws.Range("J58").Value = 1496
$ws.Range("K58").Value = 1728
$ws.Range("L58").Value = 4488
$ws.Range("M58").Value = -1600
$ws.Range("N58").Value = -4744
$ws.Range("H93").Value = 3980
$ws.Range("J93").Value = 3980
$ws.Range("L93").Value = 11940
$ws.Range("N93").Value = -15684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 204577.5
$ws.Range("I70").Value = 336961
$ws.Range("J70").Value = 6002.25
$ws.Range("K70").Value = 336961
$ws.Range("L70").Value = 6002.25
$ws.Range("M70").Value = -336691
$ws.Range("N70").Value = -6542.25
$ws.Range("H73").Value = 204577.5
$ws.Range("I73").Value = 336961
$ws.Range("J73").Value = 6002.25
$ws.Range("K73").Value = 336961
$ws.Range("L73").Value = 6002.25
$ws.Range("M73").Value = -336025
$ws.Range("N73").Value = -7874.25
$ws.Range("H136").Value = 11242.6
$ws.Range("J136").Value = 11242.6
$ws.Range("L136").Value = 33727.8
$ws.Range("N136").Value = -38827.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 53464.844
$ws.Range("I16").Value = 83935.836
$ws.Range("J16").Value = 1228.8572
$ws.Range("K16").Value = 83935.836
$ws.Range("L16").Value = 1228.8572
$ws.Range("M16").Value = -83765.836
$ws.Range("N16").Value = -1568.8572
$ws.Range("H61").Value = 2456
$ws.Range("I61").Value = 800
$ws.Range("J61").Value = 2870
$ws.Range("K61").Value = 800
$ws.Range("L61").Value = 2870
$ws.Range("M61").Value = -598
$ws.Range("N61").Value = -3274
$ws.Range("H80").Value = 10196.866
$ws.Range("J80").Value = 10196.866
$ws.Range("L80").Value = 10196.866
$ws.Range("N80").Value = -12442.866
$ws.Range("H82").Value = 1073.6666
$ws.Range("I82").Value = 880.9091
$ws.Range("J82").Value = 1376.5714
$ws.Range("K82").Value = 880.9091
$ws.Range("L82").Value = 1376.5714
$ws.Range("M82").Value = -519.9091
$ws.Range("N82").Value = -2098.5714
$ws.Range("H83").Value = 10196.866
$ws.Range("J83").Value = 10196.866
$ws.Range("L83").Value = 30590.598
$ws.Range("N83").Value = -41822.598
$ws.Range("H85").Value = 1073.6666
$ws.Range("I85").Value = 880.9091
$ws.Range("J85").Value = 1376.5714
$ws.Range("K85").Value = 880.9091
$ws.Range("L85").Value = 1376.5714
$ws.Range("M85").Value = 367.0909
$ws.Range("N85").Value = -3872.5714
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 32339.5
$ws.Range("I88").Value = 9000
$ws.Range("J88").Value = 40119.332
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 40119.332
$ws.Range("M88").Value = -8572
$ws.Range("N88").Value = -40975.332
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 32339.5
$ws.Range("I91").Value = 9000
$ws.Range("J91").Value = 40119.332
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 40119.332
$ws.Range("M91").Value = -7518
$ws.Range("N91").Value = -43083.332
$ws.Range("H100").Value = 1818.75
$ws.Range("I100").Value = 1600
$ws.Range("J100").Value = 1833.3334
$ws.Range("K100").Value = 1600
$ws.Range("L100").Value = 1833.3334
$ws.Range("M100").Value = -1059
$ws.Range("N100").Value = -2915.3334
$ws.Range("H113").Value = 2456
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 2870
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 2870
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -7210
$ws.Range("H132").Value = 3119.4243
$ws.Range("I132").Value = 3231.4
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 9694.200000000001
$ws.Range("L132").Value = 5998.9998
$ws.Range("M132").Value = -7164.200000000001
$ws.Range("N132").Value = -11058.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5271.357
$ws.Range("I132").Value = 5492.4
$ws.Range("J132").Value = 4718.75
$ws.Range("K132").Value = 16477.2
$ws.Range("L132").Value = 14156.25
$ws.Range("M132").Value = -13947.2
$ws.Range("N132").Value = -19216.25
